# Refatorando o consolidador para modelo ETL
# Updates absenteeism data rows 2-11 with new sampled records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 1001;  B = "Emilly Nascimento";         C = "Engenharia";  D = "Consulta médica";     E = 4; F = 45084; G = 4437.27 }
    @{ Row = 3;  A = 17239; B = "Davi Alves";                 C = "TI";          D = "Doença";              E = 4; F = 45097; G = 11973.45 }
    @{ Row = 4;  A = 81085; B = "Maria Vitória Farias";       C = "Jurídico";    D = "Consulta médica";     E = 4; F = 45090; G = 4747.39 }
    @{ Row = 5;  A = 28286; B = "Brenda Nascimento";          C = "Operações";   D = "Viagem de negócios";  E = 4; F = 45102; G = 7960.48 }
    @{ Row = 6;  A = 70199; B = "Ana Sophia da Mota";         C = "Jurídico";    D = "Outros";              E = 2; F = 45106; G = 3845.91 }
    @{ Row = 7;  A = 71223; B = "Sophie da Rosa";             C = "Operações";   D = "Problemas pessoais";  E = 2; F = 45084; G = 8614.780000000001 }
    @{ Row = 8;  A = 66834; B = "Pedro Ramos";                C = "Jurídico";    D = "Problemas pessoais";  E = 6; F = 45085; G = 3824 }
    @{ Row = 9;  A = 12116; B = "Juliana Barros";             C = "Engenharia";  D = "Outros";              E = 3; F = 45089; G = 6844.04 }
    @{ Row = 10; A = 2609;  B = "Dra. Ana Beatriz da Cruz";   C = "Marketing";   D = "Consulta médica";     E = 5; F = 45090; G = 2762.07 }
    @{ Row = 11; A = 48554; B = "Anthony Monteiro";           C = "Vendas";      D = "Outros";              E = 2; F = 45092; G = 5636.02 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
